$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BuildItemData")

# Write cells in the same order the original author typed them so the
# shared-string table gets built up in a matching sequence.
# Row 3 (오두막 / Cabin) - PrefabPath then SpritePath
$ws.Range("I3").Value = "Prefabs/BuildItemPrefabs/Cabin"
$ws.Range("H3").Value = "Sprites/Cabin"

# Row 4 (침대 / Bed) - SpritePath then PrefabPath
$ws.Range("H4").Value = "Sprites/Bed"
$ws.Range("I4").Value = "Prefabs/BuildItemPrefabs/Bed"

# Row 2 (작업대 / Table) - PrefabPath then SpritePath
$ws.Range("I2").Value = "Prefabs/BuildItemPrefabs/Table"
$ws.Range("H2").Value = "Sprites/Table"

# Resize SpritePath / PrefabPath columns to fit the new (longer) values.
$ws.Columns.Item(8).ColumnWidth = 11.571428571428571
$ws.Columns.Item(9).ColumnWidth = 26.285714285714285

$ws.Range("H4").Select() | Out-Null
